# ------------------------------------------------------------------
# Edit: rename "Requested quantity" header columns + add "PO Forecast"
# sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Rename headers on the existing sheets -----------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Match the outline settings used on the other sheets (summary rows
# below detail, summary columns to the right).
$ws3.Outline.SummaryRow    = 1
$ws3.Outline.SummaryColumn = 1

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Match formatting of the existing sheets: bold/centered/bordered header,
# and date-number-formatted first column.
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Data rows
$ws3.Cells.Item(2,1).Value = 45340.99999999999
$ws3.Cells.Item(2,2).Value = 33
$ws3.Cells.Item(2,3).Value = -0.7692069887035649
$ws3.Cells.Item(2,4).Value = 69.0520796563995
$ws3.Cells.Item(3,1).Value = 45347.99999999999
$ws3.Cells.Item(3,2).Value = 34
$ws3.Cells.Item(3,3).Value = -7.955532883920144
$ws3.Cells.Item(3,4).Value = 67.44523191730671
$ws3.Cells.Item(4,1).Value = 45354.99999999999
$ws3.Cells.Item(4,2).Value = 34
$ws3.Cells.Item(4,3).Value = -2.338065547246209
$ws3.Cells.Item(4,4).Value = 69.26757157287912
$ws3.Cells.Item(5,1).Value = 45361.99999999999
$ws3.Cells.Item(5,2).Value = 34
$ws3.Cells.Item(5,3).Value = -3.204035835599543
$ws3.Cells.Item(5,4).Value = 72.79377594208064
$ws3.Cells.Item(6,1).Value = 45375.99999999999
$ws3.Cells.Item(6,2).Value = 34
$ws3.Cells.Item(6,3).Value = -2.010538309984279
$ws3.Cells.Item(6,4).Value = 72.62217035021941
$ws3.Cells.Item(7,1).Value = 45389.99999999999
$ws3.Cells.Item(7,2).Value = 34
$ws3.Cells.Item(7,3).Value = -1.890749618211533
$ws3.Cells.Item(7,4).Value = 72.48492553255309
$ws3.Cells.Item(8,1).Value = 45403.99999999999
$ws3.Cells.Item(8,2).Value = 34
$ws3.Cells.Item(8,3).Value = -3.246158272610398
$ws3.Cells.Item(8,4).Value = 69.88030073555178
$ws3.Cells.Item(9,1).Value = 45431.99999999999
$ws3.Cells.Item(9,2).Value = 34
$ws3.Cells.Item(9,3).Value = -4.083461251317643
$ws3.Cells.Item(9,4).Value = 73.41697274622136
$ws3.Cells.Item(10,1).Value = 45445.99999999999
$ws3.Cells.Item(10,2).Value = 34
$ws3.Cells.Item(10,3).Value = -2.743391715239746
$ws3.Cells.Item(10,4).Value = 70.73852926483931
$ws3.Cells.Item(11,1).Value = 45452.99999999999
$ws3.Cells.Item(11,2).Value = 34
$ws3.Cells.Item(11,3).Value = -3.528125926273859
$ws3.Cells.Item(11,4).Value = 72.98435966249788
$ws3.Cells.Item(12,1).Value = 45459.99999999999
$ws3.Cells.Item(12,2).Value = 35
$ws3.Cells.Item(12,3).Value = -1.216503804430944
$ws3.Cells.Item(12,4).Value = 71.8259297486007
$ws3.Cells.Item(13,1).Value = 45473.99999999999
$ws3.Cells.Item(13,2).Value = 35
$ws3.Cells.Item(13,3).Value = -4.044961993805395
$ws3.Cells.Item(13,4).Value = 70.68043647529622
$ws3.Cells.Item(14,1).Value = 45480.99999999999
$ws3.Cells.Item(14,2).Value = 35
$ws3.Cells.Item(14,3).Value = -3.658791758507211
$ws3.Cells.Item(14,4).Value = 71.62868853818114
$ws3.Cells.Item(15,1).Value = 45487.99999999999
$ws3.Cells.Item(15,2).Value = 35
$ws3.Cells.Item(15,3).Value = -0.2416966613679906
$ws3.Cells.Item(15,4).Value = 70.98745989549245
$ws3.Cells.Item(16,1).Value = 45494.99999999999
$ws3.Cells.Item(16,2).Value = 35
$ws3.Cells.Item(16,3).Value = -3.333810700980243
$ws3.Cells.Item(16,4).Value = 71.25297659124975
$ws3.Cells.Item(17,1).Value = 45501.99999999999
$ws3.Cells.Item(17,2).Value = 35
$ws3.Cells.Item(17,3).Value = -2.970782527550342
$ws3.Cells.Item(17,4).Value = 70.67009716549551
$ws3.Cells.Item(18,1).Value = 45508.99999999999
$ws3.Cells.Item(18,2).Value = 35
$ws3.Cells.Item(18,3).Value = -1.151534300639297
$ws3.Cells.Item(18,4).Value = 74.04969728370642
$ws3.Cells.Item(19,1).Value = 45515.99999999999
$ws3.Cells.Item(19,2).Value = 35
$ws3.Cells.Item(19,3).Value = -0.8246183840761592
$ws3.Cells.Item(19,4).Value = 73.9636143809456
$ws3.Cells.Item(20,1).Value = 45522.99999999999
$ws3.Cells.Item(20,2).Value = 35
$ws3.Cells.Item(20,3).Value = -3.082794757272273
$ws3.Cells.Item(20,4).Value = 70.72717948488932
$ws3.Cells.Item(21,1).Value = 45529.99999999999
$ws3.Cells.Item(21,2).Value = 35
$ws3.Cells.Item(21,3).Value = -0.2571064366668978
$ws3.Cells.Item(21,4).Value = 70.88579969423517

$null = $ws3.Range("A1").Select()
